# day 63, 64, 65
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Fixed "date and time" placeholder text: 3/22/2018 -> 3/28/2018
#    These live on the handout master, the slide master, and every slide
#    layout (not on the slides themselves).
# ---------------------------------------------------------------------------
function Update-DateShapes($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $sh = $container.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "3/22/2018") {
                $tr.Text = "3/28/2018"
            }
        }
    }
}

Update-DateShapes $p.HandoutMaster
Update-DateShapes $p.SlideMaster

$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateShapes $layouts.Item($li)
}

# ---------------------------------------------------------------------------
# 2) Slide 3 ("Objective/To Do for Today") content placeholder edits.
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$body = $s3.Shapes.Item("Content Placeholder 2")
$tr = $body.TextFrame.TextRange

# "Essentially 13 class periods left to complete " -> "...12..."
# split the run so it matches "Essentially " / "12 " / "class periods..."
$para2 = $tr.Paragraphs(2)
$digits = $para2.Characters(13, 3)
$digits.Text = "12 "

# "Continue to " + "build" -> single run "Continue to build"
$para3 = $tr.Paragraphs(3)
$wholePara3 = $para3.Characters(1, $para3.Length)
$wholePara3.Text = "Continue to build"
